$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Quantum Computing: Unveiling Mysteries", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Unveiling the Enigma of Chemistry: The Symphony of Elements", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Author line ("Dr. Ethan Carter" -> "Amelia Stevens")
# ------------------------------------------------------------------
$pAuthor = $d.Paragraphs.Item(2)
$rAuthor = $pAuthor.Range
$rAuthor.MoveEnd(1, -1) | Out-Null
$rAuthor.Text = "Amelia Stevens"

# ------------------------------------------------------------------
# 3. Email line ("Ethan.Carter@QuantTech.Org" -> "ameliasievans@emailworld.com")
# ------------------------------------------------------------------
$pEmail = $d.Paragraphs.Item(3)
$rEmail = $pEmail.Range
$rEmail.MoveEnd(1, -1) | Out-Null
$rEmail.Text = "ameliasievans@emailworld.com"

# ------------------------------------------------------------------
# 4. Main body paragraph
# ------------------------------------------------------------------
$vb = [char]11

$bodyText = "The world of Chemistry is a fascinating tapestry of elements and molecules, where intricate interactions orchestrate the symphony of life" + `
    "." + `
    " This realm of science unveils the enigmatic secrets of matter, revealing the fundamental building blocks of our universe and the processes that shape our existence" + `
    "." + `
    " From the smallest atom to the vast expanse of galaxies, Chemistry holds the key to comprehending the delicate balance and beauty of the cosmos" + `
    "." + `
    $vb + $vb + "As we delve into the depths of Chemistry, we unravel the captivating artistry of chemical reactions, where elements dance in a harmonious ballet, transforming into new substances with novel properties" + `
    "." + `
    " The study of Chemistry empowers us to understand the intricate web of life, from the intricate workings of cellular processes to the vast array of compounds that make up the natural world" + `
    "." + `
    " It grants us the tools to unravel the mysteries of disease, paving the way for innovative treatments and therapies" + `
    "." + `
    $vb + $vb + "Furthermore, Chemistry plays a pivotal role in addressing global challenges, such as the development of sustainable energy sources, the creation of innovative materials, and the quest for cleaner and safer technologies" + `
    "." + `
    " By harnessing the power of Chemistry, we can create a sustainable future, ensuring the well-being of generations to come" + `
    "."

$pBody = $d.Paragraphs.Item(5)
$rBody = $pBody.Range
$rBody.MoveEnd(1, -1) | Out-Null
$rBody.Text = $bodyText

# ------------------------------------------------------------------
# 5. Summary paragraph (the text under the "Summary" heading)
# ------------------------------------------------------------------
$summaryText = "In conclusion, Chemistry is a captivating science that unveils the mysteries of matter, unravels the intricate symphony of chemical reactions, and empowers us to address global challenges" + `
    "." + `
    " Its study provides a profound understanding of the universe, allowing us to harness the power of elements and molecules to create innovative solutions and shape a better world" + `
    "."

$pSummary = $d.Paragraphs.Item(7)
$rSummary = $pSummary.Range
$rSummary.MoveEnd(1, -1) | Out-Null
$rSummary.Text = $summaryText

# ------------------------------------------------------------------
# 6. Append a new empty paragraph at the very end of the document
# ------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0) | Out-Null
$endRange.InsertParagraphAfter() | Out-Null

# ------------------------------------------------------------------
# 7. Fix the font everywhere: TimesNewToman -> Times New Roman
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Name = "Times New Roman"
}

Write-Host "edit complete"
